$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.815.27'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.99%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.990.05'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.33%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.82'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +9.33%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.74'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +7.27%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.684'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.55%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.749'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.96%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.168'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.73%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.71'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.55%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000318'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.46%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.98'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.629.56'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.55%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.005.53'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.83%  '

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +8.54%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.13'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.39'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.44%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.12%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.621.50'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.76%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '433.75'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.39%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.79'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +13.42%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '96.26'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.36%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.43'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.46%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.20'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.03%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.42'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +18.57%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.19%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.96'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.95%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.47'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.75%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.39'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.77%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.88'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.80'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.48%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.131'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.11%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '48.79'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '671.26'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.41%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '70.71'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.89%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.439'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0874'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.71%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.11%  '

$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.36'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.01%  '

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Dai'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.55%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.33'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.31%  '

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.36%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0490'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.63%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.72'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +9.88%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.65%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.45'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.21%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.60'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.33%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.855.23'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.57%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.03'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.70%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.02%  '
